$d = $word.ActiveDocument

# Namespace prefix used for fragment XML inserted via Range.InsertXML
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# An "empty" checklist paragraph (same formatting used throughout this
# document's checkbox column).
$emptyPara = '<w:p ' + $wns + '>' + `
  '<w:pPr>' + `
    '<w:widowControl w:val="0"/>' + `
    '<w:autoSpaceDE w:val="0"/>' + `
    '<w:autoSpaceDN w:val="0"/>' + `
    '<w:adjustRightInd w:val="0"/>' + `
    '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
'</w:p>'

# A paragraph carrying the Wingdings "checked box" glyph used to mark an
# evaluation item as complete.
$checkPara = '<w:p ' + $wns + '>' + `
  '<w:pPr>' + `
    '<w:widowControl w:val="0"/>' + `
    '<w:autoSpaceDE w:val="0"/>' + `
    '<w:autoSpaceDN w:val="0"/>' + `
    '<w:adjustRightInd w:val="0"/>' + `
    '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/>' + `
    '</w:rPr>' + `
    '<w:sym w:font="Wingdings" w:char="F0FC"/>' + `
  '</w:r>' + `
'</w:p>'

# Locate the "Problem 1 - Product Classes" evaluation table: the row whose
# first column lists "Source code includes xml comments..." and whose
# second (checkbox) column currently holds four empty paragraphs.
$targetTable = $null
$targetRow = 0
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Columns.Count -ne 2) { continue }
    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        $c1 = $null
        try { $c1 = $t.Cell($r, 1) } catch { continue }
        if ($c1.Range.Text -like "*Source code includes xml comments*" -and `
            $c1.Range.Text -like "*ProductSQLDB*") {
            $targetTable = $t
            $targetRow = $r
            break
        }
    }
    if ($targetTable -ne $null) { break }
}

$cell = $targetTable.Cell($targetRow, 2)

# Insert the seven new paragraphs (six blank + one checkmark) in front of
# the existing (empty) paragraphs already in the cell.
$firstPara = $cell.Range.Paragraphs.Item(1)
$insertRange = $firstPara.Range
$insertRange.Collapse(1)

$chunk = $emptyPara + $checkPara + $emptyPara + $emptyPara + $emptyPara + $emptyPara + $emptyPara
[void]$insertRange.InsertXML($chunk)
